# #5: property boat&car done
#
# The "汽車" (car) sheet's header row (row 1) had erroneously been filled
# with a duplicate of the data row instead of real column labels, and the
# data row (row 2) was missing the trailing metadata columns (property
# category / category / date / legislator name / legislator id /
# source file / index) that the "土地" (land) sheet already has.
#
# This fixes the header row to the real field names and fills in the
# missing columns H:N on the data row so both sheets share the same
# column layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Row 1: replace the bogus header values with the real column names ---
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# the new header cells (H1:N1) need the same bold/bordered header style
# that B1:G1 already carry
$ws.Range("B1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 2: fill in the remaining metadata columns for the car record ---
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"

# force this one to stay plain text instead of being auto-parsed as a date
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2012-04-26"

$ws.Range("K2").Value = "孔文吉"
$ws.Range("L2").Value = 1312
$ws.Range("M2").Value = "tmpfed71"
$ws.Range("N2").Value = 34

# carry the data-row style from the existing cells onto the new ones
$ws.Range("B2").Copy()
$ws.Range("H2:N2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
